# Applies:
#  1) tableStyleId change on the tables in slides 14, 15 and 16
#     ({8AB6DFCF-D2EE-4989-8ED9-20A82FBF2AFE} -> {21BFE9BD-F9E3-4290-842B-EFB92A566066})
#  2) the theme "swap" between theme1.xml (Office Theme / Office colours) and
#     theme2.xml (Integral / Red Violet colours) -- after the edit the colour
#     scheme actually used by the (sole) slide master becomes the plain
#     "Office" palette.

$p = $ppt.ActivePresentation

# -- 1. Table styles -------------------------------------------------------
$newStyleId = "{21BFE9BD-F9E3-4290-842B-EFB92A566066}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newStyleId)
        }
    }
}

# -- 2. Theme colour scheme -------------------------------------------------
# Convert a "RRGGBB" hex string into the little-endian COLORREF integer that
# the PowerPoint object model's RGB property expects (0x00BBGGRR).
function ConvertTo-Colorref([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $b * 65536 + $g * 256 + $r
}

$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$master = $p.SlideMaster
$scheme = $master.ColorScheme
for ($i = 1; $i -le $officeColors.Count; $i++) {
    $scheme.Colors($i).RGB = ConvertTo-Colorref $officeColors[$i - 1]
}
